$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ZAMORA TAMAY NEYSER IVAN", 101),
    @("MEDINA VALLEJOS ERICK LEONARDO", 98),
    @("SOTO VALLEJOS ELSITA", 96),
    @("VASQUEZ DIAZ LUZ ANGELICA", 96),
    @("TIRADO PEREZ JEINER", 95),
    @("ROJAS VASQUEZ FLOR NOELITA", 95),
    @("SOTO VILLENA NILSON", 94),
    @("BENAVIDES MARRUFO ARACELI", 94),
    @("PÓSITO CHUGDEN NANIX", 92),
    @("TELLO FERNANDEZ MILENY", 91),
    @("VASQUEZ LUNA YUDITH", 87),
    @("GALLARDO CORTEZ MELISSA DEL CARMEN", 87),
    @("BENAVIDES SALAZAR IDELSA", 81),
    @("RUIZ RUIZ LUZ MERI", 69)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
